$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 49-50, pushing the existing data (old rows 49..151)
# down to rows 51..153. Excel carries the formatting (e.g. the date number
# format in column D) from the row above when inserting, matching the
# target workbook's s="2" style on the new D49/D50 cells.
$ws.Rows("49:50").Insert()

# New row 49: Zafiro rojo, Primera, Region de Arica y Parinacota
$ws.Cells.Item(49, 1).Value = 7
$ws.Cells.Item(49, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(49, 3).Value = "Ñuble"
$ws.Cells.Item(49, 4).Value = 44477
$ws.Cells.Item(49, 5).Value = 16
$ws.Cells.Item(49, 6).Value = 100112002
$ws.Cells.Item(49, 7).Value = "Pimiento"
$ws.Cells.Item(49, 8).Value = "Zafiro rojo"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 60
$ws.Cells.Item(49, 11).Value = 43000
$ws.Cells.Item(49, 12).Value = 44000
$ws.Cells.Item(49, 13).Value = 43500
$ws.Cells.Item(49, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(49, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value = 2900
$ws.Cells.Item(49, 17).Value = 15
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# New row 50: Zafiro verde, Primera, Region de Arica y Parinacota
$ws.Cells.Item(50, 1).Value = 7
$ws.Cells.Item(50, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(50, 3).Value = "Ñuble"
$ws.Cells.Item(50, 4).Value = 44477
$ws.Cells.Item(50, 5).Value = 16
$ws.Cells.Item(50, 6).Value = 100112002
$ws.Cells.Item(50, 7).Value = "Pimiento"
$ws.Cells.Item(50, 8).Value = "Zafiro verde"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 120
$ws.Cells.Item(50, 11).Value = 41000
$ws.Cells.Item(50, 12).Value = 42000
$ws.Cells.Item(50, 13).Value = 41500
$ws.Cells.Item(50, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(50, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value = 2767
$ws.Cells.Item(50, 17).Value = 15
$ws.Cells.Item(50, 18).Value = "Hortaliza"
